$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new blank rows above row 800, shifting the old rows 800-901
# down to 802-903.
$ws.Rows.Item(800).Resize(2).Insert()

# Populate new row 800 (carries the same "shape" as the old row 800 it
# displaced, with the following fields updated per the edit).
$ws.Cells.Item(800, 1).Value = 11
$ws.Cells.Item(800, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(800, 3).Value = "Bíobío"
$ws.Cells.Item(800, 4).Value = 45154
$ws.Cells.Item(800, 5).Value = 8
$ws.Cells.Item(800, 6).Value = 100112004
$ws.Cells.Item(800, 7).Value = "Cebolla"
$ws.Cells.Item(800, 8).Value = "Sin especificar"
$ws.Cells.Item(800, 9).Value = "1a (guarda)"
$ws.Cells.Item(800, 10).Value = 200
$ws.Cells.Item(800, 11).Value = 11000
$ws.Cells.Item(800, 12).Value = 11000
$ws.Cells.Item(800, 13).Value = 11000
$ws.Cells.Item(800, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(800, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(800, 16).Value = 611
$ws.Cells.Item(800, 17).Value = 18
$ws.Cells.Item(800, 18).Value = "Hortaliza"

# Populate new row 801.
$ws.Cells.Item(801, 1).Value = 11
$ws.Cells.Item(801, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(801, 3).Value = "Bíobío"
$ws.Cells.Item(801, 4).Value = 45154
$ws.Cells.Item(801, 5).Value = 8
$ws.Cells.Item(801, 6).Value = 100112004
$ws.Cells.Item(801, 7).Value = "Cebolla"
$ws.Cells.Item(801, 8).Value = "Sin especificar"
$ws.Cells.Item(801, 9).Value = "2a (guarda)"
$ws.Cells.Item(801, 10).Value = 150
$ws.Cells.Item(801, 11).Value = 9000
$ws.Cells.Item(801, 12).Value = 9000
$ws.Cells.Item(801, 13).Value = 9000
$ws.Cells.Item(801, 14).Value = "`$/malla 18 kilos"
$ws.Cells.Item(801, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(801, 16).Value = 500
$ws.Cells.Item(801, 17).Value = 18
$ws.Cells.Item(801, 18).Value = "Hortaliza"
